$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-CellText $ws.Range("D2") "26.747.92"
Set-CellText $ws.Range("E2") "  -1.17%  "

# Row 3
Set-CellText $ws.Range("D3") "1.805.52"
Set-CellText $ws.Range("E3") "  +0.36%  "

# Row 4
Set-CellText $ws.Range("D4") "1.001"
Set-CellText $ws.Range("E4") "  -0.13%  "

# Row 5
Set-CellText $ws.Range("B5") "USDC"
Set-CellText $ws.Range("C5") "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-CellText $ws.Range("D5") "1.000"
Set-CellText $ws.Range("E5") "  -0.17%  "

# Row 6
Set-CellText $ws.Range("B6") "BNB"
Set-CellText $ws.Range("C6") "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-CellText $ws.Range("D6") "305.57"
Set-CellText $ws.Range("E6") "  -0.73%  "

# Row 7
Set-CellText $ws.Range("D7") "0.4303"
Set-CellText $ws.Range("E7") "  +2.25%  "

# Row 8
Set-CellText $ws.Range("D8") "0.3661"
Set-CellText $ws.Range("E8") "  +1.93%  "

# Row 9
Set-CellText $ws.Range("D9") "0.07211"
Set-CellText $ws.Range("E9") "  +1.42%  "

# Row 10
Set-CellText $ws.Range("D10") "0.8600"
Set-CellText $ws.Range("E10") "  +1.90%  "

# Row 11
Set-CellText $ws.Range("D11") "20.75"
Set-CellText $ws.Range("E11") "  +2.78%  "

# Row 12
Set-CellText $ws.Range("D12") "1.930.44"
Set-CellText $ws.Range("E12") "  +7.16%  "

# Row 13
Set-CellText $ws.Range("D13") "6.583"
Set-CellText $ws.Range("E13") "  +3.38%  "

# Row 14
Set-CellText $ws.Range("D14") "5.326"
Set-CellText $ws.Range("E14") "  +0.62%  "

# Row 15
Set-CellText $ws.Range("D15") "0.06892"
Set-CellText $ws.Range("E15") "  +1.86%  "

# Row 16
Set-CellText $ws.Range("D16") "1.006"
Set-CellText $ws.Range("E16") "  +0.35%  "

# Row 17
Set-CellText $ws.Range("D17") "80.22"
Set-CellText $ws.Range("E17") "  +0.19%  "

# Row 18
Set-CellText $ws.Range("D18") "0.000008878"
Set-CellText $ws.Range("E18") "  +1.77%  "

# Row 19
Set-CellText $ws.Range("D19") "1.0000"
Set-CellText $ws.Range("E19") "  -0.17%  "

# Row 20
Set-CellText $ws.Range("D20") "15.22"
Set-CellText $ws.Range("E20") "  +1.42%  "

# Row 21
Set-CellText $ws.Range("D21") "26.788.17"
Set-CellText $ws.Range("E21") "  -1.01%  "

# Row 22
Set-CellText $ws.Range("D22") "5.181"
Set-CellText $ws.Range("E22") "  +2.36%  "

# Row 23
Set-CellText $ws.Range("D23") "11.10"
Set-CellText $ws.Range("E23") "  +0.80%  "

# Row 24
Set-CellText $ws.Range("D24") "2.145.48"
Set-CellText $ws.Range("E24") "  +6.37%  "

# Row 25
Set-CellText $ws.Range("D25") "152.53"
Set-CellText $ws.Range("E25") "  -0.24%  "

# Row 26
Set-CellText $ws.Range("D26") "1.854"
Set-CellText $ws.Range("E26") "  -3.68%  "

# Row 27
Set-CellText $ws.Range("D27") "18.24"
Set-CellText $ws.Range("E27") "  +0.60%  "

# Row 28
Set-CellText $ws.Range("D28") "5.210"
Set-CellText $ws.Range("E28") "  +3.90%  "

# Row 29
Set-CellText $ws.Range("D29") "1.894"
Set-CellText $ws.Range("E29") "  +15.05%  "

# Row 30
Set-CellText $ws.Range("D30") "115.21"
Set-CellText $ws.Range("E30") "  +1.69%  "

# Row 31
Set-CellText $ws.Range("D31") "0.08923"
Set-CellText $ws.Range("E31") "  -0.69%  "

# Row 32
Set-CellText $ws.Range("D32") "0.7503"
Set-CellText $ws.Range("E32") "  +3.69%  "

# Row 33
Set-CellText $ws.Range("D33") "1.160"
Set-CellText $ws.Range("E33") "  +6.72%  "

# Row 34
Set-CellText $ws.Range("D34") "4.399"
Set-CellText $ws.Range("E34") "  +1.61%  "

# Row 35
Set-CellText $ws.Range("D35") "2.763"
Set-CellText $ws.Range("E35") "  -3.74%  "

# Row 36
Set-CellText $ws.Range("D36") "1.003"
Set-CellText $ws.Range("E36") "  +0.12%  "

# Row 37
Set-CellText $ws.Range("D37") "1.130"
Set-CellText $ws.Range("E37") "  +4.56%  "

# Row 38
Set-CellText $ws.Range("D38") "0.05184"
Set-CellText $ws.Range("E38") "  +0.95%  "

# Row 39
Set-CellText $ws.Range("D39") "0.01911"
Set-CellText $ws.Range("E39") "  +0.33%  "

# Row 40
Set-CellText $ws.Range("D40") "0.5055"
Set-CellText $ws.Range("E40") "  +1.85%  "

# Row 41
Set-CellText $ws.Range("D41") "0.1639"
Set-CellText $ws.Range("E41") "  +0.63%  "

# Row 42
Set-CellText $ws.Range("D42") "2.631"
Set-CellText $ws.Range("E42") "  +0.73%  "

# Row 43
Set-CellText $ws.Range("D43") "6.516"
Set-CellText $ws.Range("E43") "  +10.55%  "

# Row 44
Set-CellText $ws.Range("D44") "8.275"
Set-CellText $ws.Range("E44") "  +2.98%  "

# Row 45
Set-CellText $ws.Range("B45") "PaxosStandard"
Set-CellText $ws.Range("C45") "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
Set-CellText $ws.Range("D45") "1.002"
Set-CellText $ws.Range("E45") "  -30.99%  "

# Row 46
Set-CellText $ws.Range("B46") "Quant"
Set-CellText $ws.Range("C46") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-CellText $ws.Range("D46") "106.28"
Set-CellText $ws.Range("E46") "  +1.16%  "

# Row 47
Set-CellText $ws.Range("B47") "EnergySwap"
Set-CellText $ws.Range("C47") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws.Range("D47") "10.30"
Set-CellText $ws.Range("E47") "  +0.72%  "

# Row 48
Set-CellText $ws.Range("B48") "PaxDollar"
Set-CellText $ws.Range("C48") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-CellText $ws.Range("D48") "0.9995"
Set-CellText $ws.Range("E48") "  -0.21%  "

# Row 49
Set-CellText $ws.Range("B49") "NEARProtocol"
Set-CellText $ws.Range("C49") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-CellText $ws.Range("D49") "1.645"
Set-CellText $ws.Range("E49") "  +2.87%  "

# Row 50
Set-CellText $ws.Range("B50") "Decentraland"
Set-CellText $ws.Range("C50") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-CellText $ws.Range("D50") "0.4546"
Set-CellText $ws.Range("E50") "  +0.46%  "

# Row 51
Set-CellText $ws.Range("B51") "Cronos"
Set-CellText $ws.Range("C51") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText $ws.Range("D51") "0.06245"
Set-CellText $ws.Range("E51") "  -0.80%  "
